$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update milestone times for the time log entries (rows 8-12)
$ws.Range("D8").Value = 0.79791666666666661
$ws.Range("F8").Value = 0.81805555555555554

$ws.Range("C9").Value = 43567
$ws.Range("D9").Value = 0.81805555555555554
$ws.Range("F9").Value = 0.85138888888888886

$ws.Range("C10").Value = 43567
$ws.Range("D10").Value = 0.90138888888888891
$ws.Range("F10").Value = 0.95000000000000007

$ws.Range("C11").Value = 43567
$ws.Range("D11").Value = 0.95624999999999993
$ws.Range("F11").Value = 0.97986111111111107

$ws.Range("C12").Value = 43567
$ws.Range("D12").Value = 0.97986111111111107
$ws.Range("F12").Value = 0.99930555555555556

# Update selection to F12 on Sheet1
$ws.Range("F12").Select()
